# Split the run " на книга" (lang=bg-BG) into " на " (unchanged) and a new
# bold run "книга" (b + bCs), matching the commit's formatting change.

$d = $word.ActiveDocument

# Locate the target word "книга" inside the " на книга" run.
$target = $d.Content
$targetFound = $target.Find.Execute("книга", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($targetFound) {
    # Borrow formatting (bold + bold-complex-script) from an existing bold
    # run elsewhere in the same paragraph ("публикуване") so that both
    # <w:b/> and <w:bCs/> end up on the new run, then restore the original
    # "книга" text (FormattedText assignment copies text + formatting from
    # the source range).
    $src = $d.Content
    $srcFound = $src.Find.Execute("публикуване", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($srcFound) {
        $srcLen = $src.Text.Length

        $target.FormattedText = $src.FormattedText
        $target.End = $target.Start + $srcLen
        $target.Text = "книга"
    }
}
